$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.351.93'
$ws.Range('E2').Value = '  -3.51%  '
$ws.Range('D3').Value = '3.158.21'
$ws.Range('E3').Value = '  -2.97%  '
$ws.Range('E4').Value = '  +0.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.53'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.62%  '
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').Value = '3.153.65'
$ws.Range('E8').Value = '  -3.14%  '
$ws.Range('E9').Value = '  -4.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.151'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.46'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -8.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.477'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000251'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -7.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.55'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -9.64%  '
$ws.Range('D15').Value = '3.686.39'
$ws.Range('E15').Value = '  -2.48%  '
$ws.Range('D16').Value = '64.339.69'
$ws.Range('E16').Value = '  -3.49%  '
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').Value = '3.165.48'
$ws.Range('E18').Value = '  -2.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.93'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -6.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '480.89'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.68'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.87%  '
$ws.Range('E22').Value = '  -5.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.75'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.67'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -8.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.57'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.59%  '
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.88'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.47'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.19'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.72'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.92%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.113'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -21.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.73'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.77%  '
$ws.Range('E33').Value = '  +0.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.24'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.65%  '
$ws.Range('E35').Value = '  -4.62%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '54.52'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.70%  '
$ws.Range('B37').Value = 'Filecoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.98'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.12%  '
$ws.Range('D38').Value = '0.0₃0726'
$ws.Range('E38').Value = '  -8.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '452.44'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.68%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.95'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -11.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0396'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.43'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.72%  '
$ws.Range('E43').Value = '  -8.78%  '
$ws.Range('D44').Value = '2.842.57'
$ws.Range('E44').Value = '  -4.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.266'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -9.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.27'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.47'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.84%  '
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.30'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.61%  '
$ws.Range('E50').Value = '  -5.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '119.56'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.59%  '
